$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.863.74"
$ws.Range("E2").Value = "  -1.33%  "

$ws.Range("D3").Value = "1.893.57"
$ws.Range("E3").Value = "  -1.41%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'0.7755"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.92%  "

$ws.Range("D6").Value = "'244.84"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.15%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "'0.3149"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -3.25%  "

$ws.Range("D9").Value = "'0.07490"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.06%  "

$ws.Range("D10").Value = "'25.49"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.42%  "

$ws.Range("D11").Value = "'0.08116"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.32%  "

$ws.Range("D12").Value = "'0.7693"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.50%  "

$ws.Range("D13").Value = "'5.471"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.13%  "

$ws.Range("D14").Value = "1.879.50"
$ws.Range("E14").Value = "  -2.17%  "

$ws.Range("D15").Value = "'92.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.78%  "

$ws.Range("D16").Value = "'6.219"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.39%  "

$ws.Range("D17").Value = "29.809.27"
$ws.Range("E17").Value = "  -1.52%  "

$ws.Range("D18").Value = "'14.00"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.69%  "

$ws.Range("D19").Value = "'244.61"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.24%  "

$ws.Range("D20").Value = "'0.000007891"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.47%  "

$ws.Range("D21").Value = "'1.000"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.09%  "

$ws.Range("D22").Value = "'8.092"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.66%  "

$ws.Range("D23").Value = "2.112.72"
$ws.Range("E23").Value = "  -3.32%  "

$ws.Range("D24").Value = "'1.001"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.10%  "

$ws.Range("D25").Value = "'0.1574"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -4.70%  "

$ws.Range("D26").Value = "'9.457"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.38%  "

$ws.Range("D27").Value = "'162.88"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.01%  "

$ws.Range("D28").Value = "'18.83"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.04%  "

$ws.Range("E29").Value = "  -5.12%  "

$ws.Range("D30").Value = "'1.430"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +2.62%  "

$ws.Range("D31").Value = "'1.552"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.00%  "

$ws.Range("D32").Value = "'4.508"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +2.89%  "

$ws.Range("D33").Value = "'4.102"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -1.11%  "

$ws.Range("D34").Value = "'0.05515"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -3.90%  "

$ws.Range("D35").Value = "'1.259"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.03%  "

$ws.Range("D36").Value = "'0.7581"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +1.15%  "

$ws.Range("D37").Value = "'0.9990"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.82%  "

$ws.Range("D38").Value = "'2.645"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -3.15%  "

$ws.Range("D39").Value = "'0.01929"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.65%  "

$ws.Range("D40").Value = "'2.789"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -1.15%  "

$ws.Range("D41").Value = "1.167.11"
$ws.Range("E41").Value = "  +12.81%  "

$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D42").Value = "'0.4456"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.95%  "

$ws.Range("B43").Value = "Aave"
$ws.Range("C43").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D43").Value = "'73.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.47%  "

$ws.Range("D44").Value = "'5.950"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.94%  "

$ws.Range("D45").Value = "'0.8476"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.75%  "

$ws.Range("B46").Value = "RenderToken"
$ws.Range("C46").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D46").Value = "'1.908"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.15%  "

$ws.Range("B47").Value = "PaxDollar"
$ws.Range("C47").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D47").Value = "'1.000"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -0.01%  "

$ws.Range("B48").Value = "Quant"
$ws.Range("C48").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D48").Value = "'102.21"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.31%  "

$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").Value = "'9.944"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.35%  "

$ws.Range("B50").Value = "SynthetixNetwork"
$ws.Range("C50").Value = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
$ws.Range("D50").Value = "'3.100"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +0.24%  "

$ws.Range("D51").Value = "'7.565"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.85%  "
